$d = $word.ActiveDocument
$vt = [char]11

# 1. Update the run timestamp in the first Heading1 paragraph.
$d.Content.Find.Execute("2025-03-01 17:55:20.667896", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-11 08:12:31.305732", 2) | Out-Null

# Item 1: gir1.2-javascriptcoregtk-4.0/stable-security
$old = "1. gir1.2-javascriptcoregtk-4.0/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(1)" + $vt + "^p" + "Name of package: gir1.2-javascriptcoregtk-4.0/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 1: $pkgname" }

# Item 2: gir1.2-javascriptcoregtk-4.1/stable-security
$old = "2. gir1.2-javascriptcoregtk-4.1/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(2)" + $vt + "^p" + "Name of package: gir1.2-javascriptcoregtk-4.1/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 2: $pkgname" }

# Item 3: gir1.2-webkit2-4.0/stable-security
$old = "3. gir1.2-webkit2-4.0/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(3)" + $vt + "^p" + "Name of package: gir1.2-webkit2-4.0/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 3: $pkgname" }

# Item 4: gir1.2-webkit2-4.1/stable-security
$old = "4. gir1.2-webkit2-4.1/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(4)" + $vt + "^p" + "Name of package: gir1.2-webkit2-4.1/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 4: $pkgname" }

# Item 5: libjavascriptcoregtk-4.0-18/stable-security
$old = "5. libjavascriptcoregtk-4.0-18/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(5)" + $vt + "^p" + "Name of package: libjavascriptcoregtk-4.0-18/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 5: $pkgname" }

# Item 6: libjavascriptcoregtk-4.1-0/stable-security
$old = "6. libjavascriptcoregtk-4.1-0/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(6)" + $vt + "^p" + "Name of package: libjavascriptcoregtk-4.1-0/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 6: $pkgname" }

# Item 7: libjavascriptcoregtk-6.0-1/stable-security
$old = "7. libjavascriptcoregtk-6.0-1/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(7)" + $vt + "^p" + "Name of package: libjavascriptcoregtk-6.0-1/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 7: $pkgname" }

# Item 8: libtasn1-6/stable-security
$old = "8. libtasn1-6/stable-security " + $vt + "^p" + "  Current version: 4.19.0-2" + $vt + "^p" + "  Update version: " + $vt + "^p" + "  (No affiliated CVE)" + $vt
$new = "(8)" + $vt + "^p" + "Name of package: libtasn1-6/stable-security " + $vt + "^p" + "Current version: 4.19.0-2 " + $vt + "^p" + "Update version: 4.19.0-2+deb12u1 " + $vt + "^p" + "No affiliated CVEs found." + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 8: $pkgname" }

# Item 9: libwebkit2gtk-4.0-37/stable-security
$old = "9. libwebkit2gtk-4.0-37/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(9)" + $vt + "^p" + "Name of package: libwebkit2gtk-4.0-37/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 9: $pkgname" }

# Item 10: libwebkit2gtk-4.1-0/stable-security
$old = "10. libwebkit2gtk-4.1-0/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(10)" + $vt + "^p" + "Name of package: libwebkit2gtk-4.1-0/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 10: $pkgname" }

# Item 11: libwebkitgtk-6.0-4/stable-security
$old = "11. libwebkitgtk-6.0-4/stable-security " + $vt + "^p" + "  Current version: 2.46.5-1~deb12u1" + $vt + "^p" + "  Update version: 2.46.6-1~deb12u1" + $vt + "^p" + "  Affiliated CVES: " + $vt + "^p" + "    CVE-2013-3667 " + $vt
$new = "(11)" + $vt + "^p" + "Name of package: libwebkitgtk-6.0-4/stable-security " + $vt + "^p" + "Current version: 2.46.5-1~deb12u1 " + $vt + "^p" + "Update version: 2.46.6-1~deb12u1 " + $vt + "^p" + "Affiliated CVES:" + $vt + "^p" + "CVE-2013-3667" + $vt
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 11: $pkgname" }

# Item 12: openssh-client/stable-security
$old = "12. openssh-client/stable-security " + $vt + "^p" + "  Current version: 1:9.2p1-2+deb12u4" + $vt + "^p" + "  Update version: 1:9.2p1-2+deb12u5" + $vt + "^p" + "  (No affiliated CVE)"
$new = "(12)" + $vt + "^p" + "Name of package: openssh-client/stable-security " + $vt + "^p" + "Current version: 9.2p1-2+deb12u4 " + $vt + "^p" + "Update version: 1:9.2p1-2+deb12u5 " + $vt + "^p" + "No affiliated CVEs found."
$ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
if (-not $ok) { Write-Host "FAILED item 12: $pkgname" }
